# Auto-update Cloudflare export: two stale DNS records for uglyfeed.com
# ("demo.uglyfeed.com" and "xml.uglyfeed.com", rows 148-149) were removed
# from the source export. Deleting the entire rows shifts every row below
# up by two, which shrinks the sheet's used range from A1:O162 to A1:O160
# and removes what were previously the last two rows (161-162).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A148:A149").EntireRow.Delete() | Out-Null
